$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("data"): add column AI (28. 9. 2021) ---
$ws1.Range("AH1").Copy($ws1.Range("AI1"))
$ws1.Range("AI1").Value = "28. 9. 2021"

$ws1.Range("AI2").Value = 0.62
$ws1.Range("AI3").Value = 0.25
$ws1.Range("AI4").Value = 0.13
$ws1.Range("AI5").Value = 0.77
$ws1.Range("AI6").Value = 0.12
$ws1.Range("AI7").Value = 0.11
$ws1.Range("AI8").Value = 0.6899999999999999
$ws1.Range("AI9").Value = 0.24
$ws1.Range("AI10").Value = 0.07000000000000001
$ws1.Range("AI11").Value = 0.46
$ws1.Range("AI12").Value = 0.35
$ws1.Range("AI13").Value = 0.19
$ws1.Range("AI14").Value = 0.47
$ws1.Range("AI15").Value = 0.33
$ws1.Range("AI16").Value = 0.2
$ws1.Range("AI17").Value = 0.55
$ws1.Range("AI18").Value = 0.29
$ws1.Range("AI19").Value = 0.16
$ws1.Range("AI20").Value = 0.68
$ws1.Range("AI21").Value = 0.23
$ws1.Range("AI22").Value = 0.09
$ws1.Range("AI23").Value = 0.75
$ws1.Range("AI24").Value = 0.17
$ws1.Range("AI25").Value = 0.08
$ws1.Range("AI26").Value = 0.64
$ws1.Range("AI27").Value = 0.26
$ws1.Range("AI28").Value = 0.1
$ws1.Range("AI29").Value = 0.45
$ws1.Range("AI30").Value = 0.33
$ws1.Range("AI31").Value = 0.22
$ws1.Range("AI32").Value = 0.68
$ws1.Range("AI33").Value = 0.22
$ws1.Range("AI34").Value = 0.1
$ws1.Range("AI35").Value = 0.55
$ws1.Range("AI36").Value = 0.29
$ws1.Range("AI37").Value = 0.16
$ws1.Range("AI38").Value = 0.62
$ws1.Range("AI39").Value = 0.25
$ws1.Range("AI40").Value = 0.13
$ws1.Range("AI41").Value = 0.66
$ws1.Range("AI42").Value = 0.22
$ws1.Range("AI43").Value = 0.12
$ws1.Range("AI44").Value = 0.58
$ws1.Range("AI45").Value = 0.28
$ws1.Range("AI46").Value = 0.14
$ws1.Range("AI47").Value = 0.63
$ws1.Range("AI48").Value = 0.24
$ws1.Range("AI49").Value = 0.13
$ws1.Range("AI50").Value = 0.65
$ws1.Range("AI51").Value = 0.23
$ws1.Range("AI52").Value = 0.12
$ws1.Range("AI53").Value = 0.62
$ws1.Range("AI54").Value = 0.26
$ws1.Range("AI55").Value = 0.12
$ws1.Range("AI56").Value = 0.54
$ws1.Range("AI57").Value = 0.33
$ws1.Range("AI58").Value = 0.13

$ws1.Range("A59").Value = "Život během pandemie, Počet protektivních aktivit, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# --- Sheet 2 ("pocetR"): add column AH (28. 9. 2021) ---
$ws2.Range("AG1").Copy($ws2.Range("AH1"))
$ws2.Range("AH1").Value = "28. 9. 2021"

$ws2.Range("AH2").Value = 1855
$ws2.Range("AH3").Value = 456
$ws2.Range("AH4").Value = 678
$ws2.Range("AH5").Value = 721
$ws2.Range("AH6").Value = 238
$ws2.Range("AH7").Value = 483
$ws2.Range("AH8").Value = 862
$ws2.Range("AH9").Value = 539
$ws2.Range("AH10").Value = 801
$ws2.Range("AH11").Value = 515
$ws2.Range("AH12").Value = 317
$ws2.Range("AH13").Value = 325
$ws2.Range("AH14").Value = 1213
$ws2.Range("AH15").Value = 897
$ws2.Range("AH16").Value = 958
$ws2.Range("AH17").Value = 965
$ws2.Range("AH18").Value = 429
$ws2.Range("AH19").Value = 216
$ws2.Range("AH20").Value = 245

$ws2.Range("A21").Value = "Život během pandemie, Počet protektivních aktivit, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
$ws2.Range("AG21").Copy($ws2.Range("AH21"))

